# TASK_19, completed, add ElementField
#
# - Задачи!D20: mark the task's "начало" (start) timestamp, clearing the
#   remaining cell border under the finished row (matches D19's style).
# - Бэклog задач!B23:C23: log the new backlog entry "сделать страницы,
#   блоки и включить блоки в страницы.(бэкграунд,чат)".
# - Selection/active-tab bookkeeping: user ended up on the backlog sheet,
#   with Задачи's cursor left on D19.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Задачи: TASK_19 (row 20) got its D (start) timestamp filled in; the cell
# drops the bottom border that used to close off the table at that row
# (now matching the borderless D19 above it).
$ws1.Range("D20").Value = 41987.771527777775
$ws1.Range("D20").NumberFormat = "dd/mm/yy\ h:mm;@"
$ws1.Range("D20").Borders.Item(9).LineStyle = -4142

# Бэклог задач: append the new backlog item.
$ws2.Range("B23").Value = "сделать страницы, блоки и включить блоки в страницы.(бэкграунд,чат)"
$ws2.Range("C23").Value = 41988.625694444447
$ws2.Range("C23").NumberFormat = "m/d/yy h:mm"

# Final selection/active sheet state as saved in the workbook.
$ws1.Range("D19").Select()
$ws2.Activate()
$ws2.Range("C24").Select()
